# Update the division problems in the practice-sheet table.
# Each "A÷B=" string appears exactly once in the document, so a simple
# whole-word Find/Replace (MatchWholeWord) per pair is sufficient.
#
# NOTE on ordering: "14÷3=" is being changed to "17÷5=", and separately the
# original "17÷5=" is being changed to "16÷3=". To avoid the second
# replacement mistakenly catching the text just produced by the first, the
# "17÷5=" -> "16÷3=" replacement is executed BEFORE the "14÷3=" -> "17÷5="
# replacement.

$d = $word.ActiveDocument

$replacements = @(
    @("96÷9=", "51÷8="),
    @("17÷5=", "16÷3="),
    @("89÷7=", "36÷6="),
    @("41÷7=", "86÷2="),
    @("53÷8=", "52÷6="),
    @("79÷6=", "83÷8="),
    @("37÷5=", "55÷7="),
    @("51÷3=", "68÷6="),
    @("80÷9=", "66÷8="),
    @("48÷8=", "38÷8="),
    @("70÷4=", "91÷3="),
    @("68÷7=", "95÷8="),
    @("11÷6=", "21÷4="),
    @("20÷9=", "71÷7="),
    @("26÷5=", "13÷3="),
    @("24÷8=", "25÷9="),
    @("33÷6=", "74÷2="),
    @("15÷5=", "98÷4="),
    @("80÷7=", "21÷5="),
    @("35÷9=", "49÷9="),
    @("14÷9=", "71÷8="),
    @("50÷4=", "99÷3="),
    @("70÷9=", "85÷5="),
    @("82÷7=", "70÷8="),
    @("14÷3=", "17÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}
